$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 7
$ws.Range("F7").Value = 1672.1
$ws.Range("G7").Value = 1681.8
$ws.Range("H7").Value = 1660.85
$ws.Range("I7").Value = 1666.85
$ws.Range("J7").Value = 1680.9

# Row 9
$ws.Range("G9").Value = 1678.45
$ws.Range("H9").Value = 1660.25
$ws.Range("I9").Value = 1671.15

# Row 10
$ws.Range("G10").Value = 1672.9
$ws.Range("H10").Value = 1662.7
$ws.Range("I10").Value = 1669.9

# Row 11
$ws.Range("G11").Value = 1680.4
$ws.Range("H11").Value = 1668.6
$ws.Range("I11").Value = 1674

# Row 12
$ws.Range("G12").Value = 1680
$ws.Range("H12").Value = 1669.8
$ws.Range("I12").Value = 1679.5

# Row 13
$ws.Range("G13").Value = 1681.8
$ws.Range("H13").Value = 1672.15
$ws.Range("I13").Value = 1676.1

# Row 14
$ws.Range("G14").Value = 1679.85
$ws.Range("H14").Value = 1675
$ws.Range("I14").Value = 1675.95

# Row 15
$ws.Range("G15").Value = 1679.95
$ws.Range("H15").Value = 1675.4
$ws.Range("I15").Value = 1678.9

# Row 16
$ws.Range("G16").Value = 1679.55
$ws.Range("H16").Value = 1674.35
$ws.Range("I16").Value = 1677.2

# Row 17
$ws.Range("G17").Value = 1679.7
$ws.Range("H17").Value = 1675
$ws.Range("I17").Value = 1678.95

# Row 18
$ws.Range("G18").Value = 1680
$ws.Range("H18").Value = 1667
$ws.Range("I18").Value = 1669.7

# Row 19
$ws.Range("G19").Value = 1671.9
$ws.Range("H19").Value = 1661
$ws.Range("I19").Value = 1668.95

# Row 20
$ws.Range("G20").Value = 1671.55
$ws.Range("H20").Value = 1660.85
$ws.Range("I20").Value = 1663.75

# Row 21
$ws.Range("G21").Value = 1672
$ws.Range("H21").Value = 1661.4
$ws.Range("I21").Value = 1671.95
